$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4 (Leve Item ID 5470)
$ws.Range("H4").Value = 92
$ws.Range("I4").Value = 92
$ws.Range("K4").Value = 92
$ws.Range("M4").Value = 22
# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 5176.6113
$ws.Range("I9").Value = 6559.357
$ws.Range("J9").Value = 337
$ws.Range("K9").Value = 6559.357
$ws.Range("L9").Value = 337
$ws.Range("M9").Value = -6390.357
$ws.Range("N9").Value = -675
# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 1270.238
$ws.Range("I28").Value = 1117.7693
$ws.Range("J28").Value = 1518
$ws.Range("K28").Value = 1117.7693
$ws.Range("L28").Value = 1518
$ws.Range("M28").Value = -632.7692999999999
$ws.Range("N28").Value = -2488
# Row 49 (Leve Item ID 4588)
$ws.Range("H49").Value = 918.75
$ws.Range("I49").Value = 154.8
$ws.Range("K49").Value = 464.4
$ws.Range("M49").Value = -328.4
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1325.5652
$ws.Range("I137").Value = 1284.25
$ws.Range("K137").Value = 3852.75
$ws.Range("M137").Value = -1302.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5 (Leve Item ID 5091)
$ws.Range("H5").Value = 50346.65
$ws.Range("I5").Value = 91171.37
$ws.Range("J5").Value = 449.77777
$ws.Range("K5").Value = 91171.37
$ws.Range("L5").Value = 449.77777
$ws.Range("M5").Value = -91059.37
$ws.Range("N5").Value = -673.7777699999999
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 5890.2856
$ws.Range("I32").Value = 5442.8335
$ws.Range("K32").Value = 5442.8335
$ws.Range("M32").Value = -5155.8335
# Row 43 (Leve Item ID 21715)
$ws.Range("H43").Value = 44751
$ws.Range("J43").Value = 44751
$ws.Range("L43").Value = 44751
$ws.Range("N43").Value = -45377
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 10574.643
$ws.Range("I45").Value = 15402.125
$ws.Range("J45").Value = 4138
$ws.Range("K45").Value = 15402.125
$ws.Range("L45").Value = 4138
$ws.Range("M45").Value = -15025.125
$ws.Range("N45").Value = -4892
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 2471.125
$ws.Range("I61").Value = 2062.8965
$ws.Range("J61").Value = 6417.3335
$ws.Range("K61").Value = 2062.8965
$ws.Range("L61").Value = 6417.3335
$ws.Range("M61").Value = -1850.8965
$ws.Range("N61").Value = -6841.3335
# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 1284.8334
$ws.Range("I63").Value = 1341.8
$ws.Range("K63").Value = 1341.8
$ws.Range("M63").Value = -655.8
# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 1284.8334
$ws.Range("I66").Value = 1341.8
$ws.Range("K66").Value = 6709
$ws.Range("M66").Value = -3277
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 8025.9614
$ws.Range("I74").Value = 867.4211
$ws.Range("K74").Value = 867.4211
$ws.Range("M74").Value = 6.578899999999976
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 8025.9614
$ws.Range("I77").Value = 867.4211
$ws.Range("K77").Value = 4337.1055
$ws.Range("M77").Value = 30.89450000000033
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 2699.8918
$ws.Range("I132").Value = 2176.8708
$ws.Range("K132").Value = 6530.6124
$ws.Range("M132").Value = -4000.6124
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 2471.125
$ws.Range("I136").Value = 2062.8965
$ws.Range("J136").Value = 6417.3335
$ws.Range("K136").Value = 6188.689499999999
$ws.Range("L136").Value = 19252.0005
$ws.Range("M136").Value = -3638.689499999999
$ws.Range("N136").Value = -24352.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4 (Leve Item ID 5091)
$ws.Range("H4").Value = 50346.65
$ws.Range("I4").Value = 91171.37
$ws.Range("J4").Value = 449.77777
$ws.Range("K4").Value = 91171.37
$ws.Range("L4").Value = 449.77777
$ws.Range("M4").Value = -91056.37
$ws.Range("N4").Value = -679.7777699999999
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 5188.1113
$ws.Range("J86").Value = 5165.8335
$ws.Range("L86").Value = 5165.8335
$ws.Range("N86").Value = -7411.8335
# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 5188.1113
$ws.Range("J89").Value = 5165.8335
$ws.Range("L89").Value = 25829.1675
$ws.Range("N89").Value = -37061.1675
# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 1442.3143
$ws.Range("I94").Value = 1209.5172
$ws.Range("J94").Value = 2567.5
$ws.Range("K94").Value = 1209.5172
$ws.Range("L94").Value = 2567.5
$ws.Range("M94").Value = -758.5172
$ws.Range("N94").Value = -3469.5
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 3460.889
$ws.Range("I105").Value = 3906.8572
$ws.Range("J105").Value = 1900
$ws.Range("K105").Value = 3906.8572
$ws.Range("L105").Value = 1900
$ws.Range("M105").Value = -2159.8572
$ws.Range("N105").Value = -5394
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 2034.3396
$ws.Range("I134").Value = 1816
$ws.Range("K134").Value = 5448
$ws.Range("M134").Value = -2913

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 35990.258
$ws.Range("J31").Value = 8232.5
$ws.Range("L31").Value = 8232.5
$ws.Range("N31").Value = -8822.5
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 35990.258
$ws.Range("J34").Value = 8232.5
$ws.Range("L34").Value = 8232.5
$ws.Range("N34").Value = -8636.5
# Row 109 (Leve Item ID 27203)
$ws.Range("H109").Value = 39333.332
$ws.Range("J109").Value = 39333.332
$ws.Range("L109").Value = 39333.332
$ws.Range("N109").Value = -41413.332

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34 (Leve Item ID 4749)
$ws.Range("H34").Value = 5556525
$ws.Range("I34").Value = 39
$ws.Range("J34").Value = 11113011
$ws.Range("K34").Value = 117
$ws.Range("L34").Value = 33339033
$ws.Range("M34").Value = -33
$ws.Range("N34").Value = -33339201
# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 1122.909
$ws.Range("J39").Value = 1225.2
$ws.Range("L39").Value = 3675.6
$ws.Range("N39").Value = -4263.6
# Row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 1074.6
$ws.Range("J55").Value = 1721.6666
$ws.Range("L55").Value = 5164.9998
$ws.Range("N55").Value = -5518.9998
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 3554.4443
$ws.Range("I122").Value = 3998.3333
$ws.Range("J122").Value = 3332.5
$ws.Range("K122").Value = 35984.9997
$ws.Range("L122").Value = 29992.5
$ws.Range("M122").Value = -33534.9997
$ws.Range("N122").Value = -34892.5
# Row 124 (Leve Item ID 36040)
$ws.Range("H124").Value = 25855.137
$ws.Range("J124").Value = 28944.475
$ws.Range("L124").Value = 86833.42499999999
$ws.Range("N124").Value = -96653.42499999999
# Row 129 (Leve Item ID 36054)
$ws.Range("H129").Value = 889.8333
$ws.Range("I129").Value = 767.8
$ws.Range("K129").Value = 2303.4
$ws.Range("M129").Value = 2696.6
# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1264.6471
$ws.Range("I132").Value = 1038.3846
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 9345.4614
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -6815.4614
$ws.Range("N132").Value = -23060
# Row 137 (Leve Item ID 44088)
$ws.Range("H137").Value = 2847.111
$ws.Range("I137").Value = 2089.2856
$ws.Range("K137").Value = 6267.8568
$ws.Range("M137").Value = -1167.8568
# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 4573.625
$ws.Range("I140").Value = 4941.4287
$ws.Range("J140").Value = 1999
$ws.Range("K140").Value = 14824.2861
$ws.Range("L140").Value = 5997
$ws.Range("M140").Value = -9644.286100000001
$ws.Range("N140").Value = -16357

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 5050
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 7600
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 7600
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -9596
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 5050
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 7600
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 38000
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -47984
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 2373.6667
$ws.Range("I113").Value = 2373.6667
$ws.Range("K113").Value = 2373.6667
$ws.Range("M113").Value = -203.6667000000002
# Row 123 (Leve Item ID 34150)
$ws.Range("H123").Value = 35700
$ws.Range("J123").Value = 35700
$ws.Range("L123").Value = 35700
$ws.Range("N123").Value = -40600
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 288065.03
$ws.Range("I132").Value = 386887.53
$ws.Range("J132").Value = 2577.7778
$ws.Range("K132").Value = 1160662.59
$ws.Range("L132").Value = 7733.3334
$ws.Range("M132").Value = -1158132.59
$ws.Range("N132").Value = -12793.3334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 11771.083
$ws.Range("I7").Value = 13027.667
$ws.Range("K7").Value = 13027.667
$ws.Range("M7").Value = -12915.667
# Row 38 (Leve Item ID 2767)
$ws.Range("H38").Value = 34000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 34000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 34000
$ws.Range("N38").Value = -34820
$ws.Range("M38").ClearContents()
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 11771.083
$ws.Range("I126").Value = 13027.667
$ws.Range("K126").Value = 39083.001
$ws.Range("M126").Value = -36613.001
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 5393
$ws.Range("I132").Value = 4570.7144
$ws.Range("J132").Value = 6112.5
$ws.Range("K132").Value = 13712.1432
$ws.Range("L132").Value = 18337.5
$ws.Range("M132").Value = -11182.1432
$ws.Range("N132").Value = -23397.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 28 (Leve Item ID 3053)
$ws.Range("H28").Value = 17381.666
$ws.Range("J28").Value = 17381.666
$ws.Range("L28").Value = 17381.666
$ws.Range("N28").Value = -18077.666
# Row 97 (Leve Item ID 18220)
$ws.Range("H97").Value = 35747
$ws.Range("J97").Value = 35747
$ws.Range("L97").Value = 35747
$ws.Range("N97").Value = -37729
# Row 109 (Leve Item ID 27161)
$ws.Range("H109").Value = 39300
$ws.Range("J109").Value = 39300
$ws.Range("L109").Value = 39300
$ws.Range("N109").Value = -42074
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 766.3333
$ws.Range("I113").Value = 749
$ws.Range("J113").Value = 775
$ws.Range("K113").Value = 2247
$ws.Range("L113").Value = 2325
$ws.Range("M113").Value = -77
$ws.Range("N113").Value = -6665
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2615.3667
$ws.Range("I132").Value = 2596.2173
$ws.Range("K132").Value = 7788.651899999999
$ws.Range("M132").Value = -5258.651899999999

Write-Host "Edit complete"
